$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values can look like numbers (e.g. "545.71"); force them to stay
# as plain text by temporarily marking the range as Text before assignment,
# then clear the temporary formatting again so no style change remains.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '56.816.93'
$ws.Range("E2").Value = '  -6.89%  '

$ws.Range("D3").Value = '2.874.94'
$ws.Range("E3").Value = '  -4.35%  '

$ws.Range("E4").Value = '  +0.34%  '

$ws.Range("D5").Value = '545.71'
$ws.Range("E5").Value = '  -3.44%  '

$ws.Range("D6").Value = '121.87'
$ws.Range("E6").Value = '  -5.09%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = '2.872.60'
$ws.Range("E8").Value = '  -4.17%  '

$ws.Range("D9").Value = '0.495'
$ws.Range("E9").Value = '  -0.49%  '

$ws.Range("D10").Value = '0.121'
$ws.Range("E10").Value = '  -10.16%  '

$ws.Range("D11").Value = '4.63'
$ws.Range("E11").Value = '  -11.24%  '

$ws.Range("D12").Value = '0.433'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").Value = '0.0000209'
$ws.Range("E13").Value = '  -6.36%  '

$ws.Range("D14").Value = '32.31'
$ws.Range("E14").Value = '  -1.86%  '

$ws.Range("E15").Value = '  +0.64%  '

$ws.Range("D16").Value = '3.349.40'
$ws.Range("E16").Value = '  -4.40%  '

$ws.Range("D17").Value = '2.881.90'
$ws.Range("E17").Value = '  -4.01%  '

$ws.Range("D18").Value = '6.49'
$ws.Range("E18").Value = '  +5.33%  '

$ws.Range("D19").Value = '57.109.67'
$ws.Range("E19").Value = '  -6.43%  '

$ws.Range("D20").Value = '400.35'
$ws.Range("E20").Value = '  -8.01%  '

$ws.Range("D21").Value = '12.76'
$ws.Range("E21").Value = '  -3.00%  '

$ws.Range("D22").Value = '0.666'
$ws.Range("E22").Value = '  +0.40%  '

$ws.Range("D23").Value = '6.78'
$ws.Range("E23").Value = '  -5.13%  '

$ws.Range("D24").Value = '12.62'
$ws.Range("E24").Value = '  -1.58%  '

$ws.Range("D25").Value = '76.73'
$ws.Range("E25").Value = '  -3.05%  '

$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("E27").Value = '  +0.73%  '

$ws.Range("D28").Value = '2.42'
$ws.Range("E28").Value = '  -2.94%  '

$ws.Range("D29").Value = '1.91'
$ws.Range("E29").Value = '  +1.54%  '

$ws.Range("D30").Value = '7.10'
$ws.Range("E30").Value = '  -1.52%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '24.50'
$ws.Range("E31").Value = '  -3.85%  '

$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = '5.91'
$ws.Range("E32").Value = '  -3.26%  '

$ws.Range("D33").Value = '0.0976'
$ws.Range("E33").Value = '  +3.93%  '

$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").Value = '0.904'
$ws.Range("E34").Value = '  -5.64%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '5.38'
$ws.Range("E35").Value = '  -3.96%  '

$ws.Range("D36").Value = '1.97'
$ws.Range("E36").Value = '  -12.94%  '

$ws.Range("D37").Value = '47.68'
$ws.Range("E37").Value = '  -4.80%  '

$ws.Range("D38").Value = '8.23'
$ws.Range("E38").Value = '  +5.78%  '

$ws.Range("D39").Value = '0.0₃0613'
$ws.Range("E39").Value = '  -8.78%  '

$ws.Range("D40").Value = '0.104'
$ws.Range("E40").Value = '  -2.98%  '

$ws.Range("D41").Value = '0.0336'
$ws.Range("E41").Value = '  -6.75%  '

$ws.Range("D42").Value = '2.611.02'
$ws.Range("E42").Value = '  -2.96%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '2.38'
$ws.Range("E43").Value = '  -2.97%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = '355.12'
$ws.Range("E44").Value = '  -6.48%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").Value = '119.04'
$ws.Range("E46").Value = '  -1.46%  '

$ws.Range("D47").Value = '0.227'
$ws.Range("E47").Value = '  -3.60%  '

$ws.Range("E48").Value = '  -0.55%  '

$ws.Range("D49").Value = '1.91'
$ws.Range("E49").Value = '  -3.06%  '

$ws.Range("D50").Value = '22.60'
$ws.Range("E50").Value = '  -3.23%  '

$ws.Range("D51").Value = '1.93'
$ws.Range("E51").Value = '  -4.87%  '

# Remove the temporary text formatting so the cells keep their original style
$dRange.ClearFormats()
